# "its roi changes for boss on maors computer"
# Reposition the two floating screenshot pictures anchored in the
# document's only paragraph, and drop the stray _GoBack bookmark that
# used to wrap the first picture's run.

$d = $word.ActiveDocument

# wdRelativeHorizontalPosition constants (not predefined in this
# stripped-down PS host, so spell them out):
#   wdRelativeHorizontalPositionMargin = 0
#   wdRelativeHorizontalPositionPage   = 1
#   wdRelativeHorizontalPositionColumn = 2
$wdRelativeHorizontalPositionPage = 1

# EMU -> points (Word's Shape.Left/Top/Width/Height are all in points;
# the OOXML stores EMUs, 12700 EMU per point).
$emuPerPt = 12700

# --- Picture 1 (relativeHeight 251658240, "image1.png") -------------
# Moves from margin-relative to page-relative, and both offsets shift.
$pic1 = $d.Shapes.Item(1)
$pic1.RelativeHorizontalPosition = $wdRelativeHorizontalPositionPage
$pic1.Left = 618646 / $emuPerPt
$pic1.Top = 4261485 / $emuPerPt
# The picture's outer box size (wp:extent) is unchanged; re-asserting
# Width/Height at that same size is what syncs the inner pic:spPr
# a:ext (previously a slightly stale 6945260x4548712) back to match.
$pic1.Width = 6941820 / $emuPerPt
$pic1.Height = 4546459 / $emuPerPt

# --- Picture 2 (relativeHeight 251659264, "image2.png") -------------
# Stays column-relative; only the horizontal offset shifts.
$pic2 = $d.Shapes.Item(2)
$pic2.Left = -993140 / $emuPerPt

# --- Drop the leftover _GoBack bookmark ------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
